$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name value in B2
$ws.Range("B2").Value = "甄任珍"

# Delete column C (Class / M120)
$ws.Range("C1:C2").Delete()
